# "Fixing Library and Fixing Excel" - update the SCD0270 test-data rows with
# new sample values and refresh the sheet's active view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (USERID / TEXT1 CIF / TEXT3 CIF)
$ws.Range("F2").Value = 35758
$ws.Range("M2").Value = 9720826341
$ws.Range("O2").Value = 9669179367
$ws.Range("O2").VerticalAlignment = -4108
# -4108 = xlCenter, matches the new style used on O2

# Row 3 (USERID only)
$ws.Range("F3").Value = 32362

# Row 4 (USERID / TEXT1 CIF / TEXT3 CIF)
$ws.Range("F4").Value = 35758
$ws.Range("M4").Value = 9720826341
$ws.Range("O4").Value = 9669179367
$ws.Range("O4").VerticalAlignment = -4108
# -4108 = xlCenter, matches the new style used on O4

# Refresh the view: scroll so column E / row 3 is top-left, select F3.
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("F3").Select()

$wb.Saved = $false
